$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I26").Value = "sd"
$ws.Range("J26").Value = "Statement-non-opinion"
$ws.Range("I32").Value = "sv"
$ws.Range("J32").Value = "Statement-opinion"
$ws.Range("I34").Value = "b"
$ws.Range("J34").Value = "Acknowledge (Backchannel)"
$ws.Range("I36").Value = "sv"
$ws.Range("J36").Value = "Statement-opinion"
$ws.Range("I44").Value = "ba"
$ws.Range("J44").Value = "Appreciation"
$ws.Range("I45").Value = "ba"
$ws.Range("J45").Value = "Appreciation"
$ws.Range("I63").Value = "sd"
$ws.Range("J63").Value = "Statement-non-opinion"
$ws.Range("I80").Value = "%"
$ws.Range("J80").Value = "Uninterpretable"
$ws.Range("I81").Value = "ba"
$ws.Range("J81").Value = "Appreciation"
$ws.Range("I83").Value = "sv"
$ws.Range("J83").Value = "Statement-opinion"
$ws.Range("I94").Value = "sd"
$ws.Range("J94").Value = "Statement-non-opinion"
$ws.Range("I106").Value = "sd"
$ws.Range("J106").Value = "Statement-non-opinion"
$ws.Range("I119").Value = "aa"
$ws.Range("J119").Value = "Agree/Accept"
$ws.Range("I128").Value = "sd"
$ws.Range("J128").Value = "Statement-non-opinion"
$ws.Range("I130").Value = "aa"
$ws.Range("J130").Value = "Agree/Accept"
$ws.Range("I133").Value = "sd"
$ws.Range("J133").Value = "Statement-non-opinion"
$ws.Range("I144").Value = "sd"
$ws.Range("J144").Value = "Statement-non-opinion"
$ws.Range("I158").Value = "sd"
$ws.Range("J158").Value = "Statement-non-opinion"
$ws.Range("I163").Value = "sd"
$ws.Range("J163").Value = "Statement-non-opinion"
$ws.Range("I172").Value = "sv"
$ws.Range("J172").Value = "Statement-opinion"
$ws.Range("I173").Value = "sd"
$ws.Range("J173").Value = "Statement-non-opinion"
$ws.Range("I188").Value = "aa"
$ws.Range("J188").Value = "Agree/Accept"
$ws.Range("I194").Value = "b"
$ws.Range("J194").Value = "Acknowledge (Backchannel)"
$ws.Range("I200").Value = "%"
$ws.Range("J200").Value = "Uninterpretable"
$ws.Range("I218").Value = "sd"
$ws.Range("J218").Value = "Statement-non-opinion"
$ws.Range("I222").Value = "ba"
$ws.Range("J222").Value = "Appreciation"
$ws.Range("I237").Value = "aa"
$ws.Range("J237").Value = "Agree/Accept"
